$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 14743
$ws1.Range("F4").Value = 18072
$ws1.Range("F5").Value = 18072
$ws1.Range("F24").Value = 7475
$ws1.Range("F30").Value = 5899

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 14743
$ws4.Range("F4").Value = 18072
$ws4.Range("F5").Value = 18072
$ws4.Range("F25").Value = 7475
$ws4.Range("F32").Value = 5899
